$d = $word.ActiveDocument

# The document contains three <id>...</id> markers whose body text used to
# carry a suffix "_aN"; they need to become a single merged run reading
# "<id>p039v_N</id>" (collapsing the three separate runs - open tag, id
# text, close tag - into one run matching the open/close tag formatting).

$d.Content.Find.Execute("<id>p039v_a1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p039v_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p039v_a2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p039v_2</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p039v_a3</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p039v_3</id>", 2) | Out-Null
